$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1.91
$ws.Range("G2").Value = 2.14
$ws.Range("H2").Value = 3.95
$ws.Range("I2").Value = 5.3
$ws.Range("J2").Value = 3.4
$ws.Range("K2").Value = 4.1
$ws.Range("P2").Value = 1.88
$ws.Range("Q2").Value = 1.92

# Row 6
$ws.Range("F6").Value = 2.44
$ws.Range("G6").Value = 2.52
$ws.Range("I6").Value = 3.5
$ws.Range("L6").Value = 1.45
$ws.Range("M6").Value = 1.1
$ws.Range("N6").Value = 3.25
$ws.Range("O6").Value = 1.42
$ws.Range("P6").Value = 1.74
$ws.Range("Q6").Value = 2.28
$ws.Range("R6").Value = 1.27
$ws.Range("U6").Value = 1.98
$ws.Range("V6").Value = 1.4
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 11.5
$ws.Range("Z6").Value = 23
$ws.Range("AB6").Value = 9
$ws.Range("AE6").Value = 44
$ws.Range("AI6").Value = 60
$ws.Range("AN6").Value = 27
$ws.Range("AO6").Value = 48

# Row 7
$ws.Range("G7").Value = 1.55
$ws.Range("J7").Value = 4.9

# Row 9
$ws.Range("K9").Value = 4

# Row 10
$ws.Range("G10").Value = 21
$ws.Range("H10").Value = 1.3

# Row 11
$ws.Range("G11").Value = 1.55
$ws.Range("I11").Value = 10
$ws.Range("P11").Value = 1.68

# Row 13
$ws.Range("J13").Value = 2.98
$ws.Range("N13").Value = 2.64
$ws.Range("O13").Value = 1.49
$ws.Range("P13").Value = 1.55
$ws.Range("Q13").Value = 2.46

# Row 14
$ws.Range("G14").Value = 1.77
$ws.Range("I14").Value = 9.2
$ws.Range("Q14").Value = 1.92

# Row 15
$ws.Range("F15").Value = 1.79
$ws.Range("I15").Value = 2.26
$ws.Range("J15").Value = 1.79
$ws.Range("N15").Value = 1.48
$ws.Range("O15").Value = 1.01

# Row 18
$ws.Range("F18").Value = 3.85
$ws.Range("T18").Value = 1.8
$ws.Range("U18").Value = 1.73
$ws.Range("X18").Value = 15
$ws.Range("Y18").Value = 10.5
$ws.Range("Z18").Value = 17
$ws.Range("AA18").Value = 38
$ws.Range("AB18").Value = 18
$ws.Range("AC18").Value = 10.5
$ws.Range("AD18").Value = 15.5
$ws.Range("AE18").Value = 36
$ws.Range("AF18").Value = 44
$ws.Range("AG18").Value = 25
$ws.Range("AH18").Value = 30

# Row 19
$ws.Range("G19").Value = 2.26
$ws.Range("AF19").Value = 13

# Row 20
$ws.Range("F20").Value = 1.88
$ws.Range("G20").Value = 1.93
$ws.Range("H20").Value = 4.2
$ws.Range("I20").Value = 4.7
$ws.Range("O20").Value = 1.28

# Row 21
$ws.Range("G21").Value = 1.89
$ws.Range("R21").Value = 1.64
$ws.Range("AE21").Value = 44
$ws.Range("AM21").Value = 55

# Row 22
$ws.Range("K22").Value = 3.6
$ws.Range("N22").Value = 3.65
$ws.Range("P22").Value = 1.91
$ws.Range("Q22").Value = 2.06

# Row 23
$ws.Range("G23").Value = 1.22
$ws.Range("H23").Value = 19
$ws.Range("J23").Value = 7.2
$ws.Range("P23").Value = 2.36
